$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Apoe"
$ws.Cells.Item(2, 3).Value = "Scarb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 29.32133366666666
$ws.Cells.Item(2, 8).Value = 87.964001
$ws.Cells.Item(2, 9).Value = 0.006401919837078288
$ws.Cells.Item(2, 10).Value = 0.006401919837078288
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 82.48638166666666
$ws.Cells.Item(2, 14).Value = 247.459145
$ws.Cells.Item(2, 15).Value = 0.7894957391680832
$ws.Cells.Item(2, 16).Value = 0.7894957391680832
$ws.Cells.Item(2, 17).Value = 2418.610719804349
$ws.Cells.Item(2, 18).Value = 21767.49647823914
$ws.Cells.Item(2, 19).Value = 0.005054288433868938
$ws.Cells.Item(2, 20).Value = 0.005054288433868939

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Apoe"
$ws.Cells.Item(3, 3).Value = "Scarb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 29.32133366666666
$ws.Cells.Item(3, 8).Value = 87.964001
$ws.Cells.Item(3, 9).Value = 0.006401919837078288
$ws.Cells.Item(3, 10).Value = 0.006401919837078288
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.8713403333333334
$ws.Cells.Item(3, 14).Value = 2.614021
$ws.Cells.Item(3, 15).Value = 0.008339794601633706
$ws.Cells.Item(3, 16).Value = 0.008339794601633706
$ws.Cells.Item(3, 17).Value = 25.54886065089122
$ws.Cells.Item(3, 18).Value = 229.939745858021
$ws.Cells.Item(3, 19).Value = 0.00005339069649735723
$ws.Cells.Item(3, 20).Value = 0.00005339069649735724

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Apoe"
$ws.Cells.Item(4, 3).Value = "Scarb1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 29.32133366666666
$ws.Cells.Item(4, 8).Value = 87.964001
$ws.Cells.Item(4, 9).Value = 0.006401919837078288
$ws.Cells.Item(4, 10).Value = 0.006401919837078288
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 17.88507033333333
$ws.Cells.Item(4, 14).Value = 53.65521099999999
$ws.Cells.Item(4, 15).Value = 0.1711820368112258
$ws.Cells.Item(4, 16).Value = 0.1711820368112258
$ws.Cells.Item(4, 17).Value = 524.4141148954678
$ws.Cells.Item(4, 18).Value = 4719.72703405921
$ws.Cells.Item(4, 19).Value = 0.001095893677213252
$ws.Cells.Item(4, 20).Value = 0.001095893677213253

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Apoe"
$ws.Cells.Item(5, 3).Value = "Scarb1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 29.32133366666666
$ws.Cells.Item(5, 8).Value = 87.964001
$ws.Cells.Item(5, 9).Value = 0.006401919837078288
$ws.Cells.Item(5, 10).Value = 0.006401919837078288
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 3.237038999999999
$ws.Cells.Item(5, 14).Value = 9.711116999999998
$ws.Cells.Item(5, 15).Value = 0.03098242941905719
$ws.Cells.Item(5, 16).Value = 0.03098242941905719
$ws.Cells.Item(5, 17).Value = 94.91430061101298
$ws.Cells.Item(5, 18).Value = 854.2287054991168
$ws.Cells.Item(5, 19).Value = 0.0001983470294987401
$ws.Cells.Item(5, 20).Value = 0.0001983470294987402

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Apoe"
$ws.Cells.Item(6, 3).Value = "Scarb1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 45.524413
$ws.Cells.Item(6, 8).Value = 136.573239
$ws.Cells.Item(6, 9).Value = 0.009939644832300594
$ws.Cells.Item(6, 10).Value = 0.009939644832300592
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 82.48638166666666
$ws.Cells.Item(6, 14).Value = 247.459145
$ws.Cells.Item(6, 15).Value = 0.7894957391680832
$ws.Cells.Item(6, 16).Value = 0.7894957391680832
$ws.Cells.Item(6, 17).Value = 3755.144105868962
$ws.Cells.Item(6, 18).Value = 33796.29695282065
$ws.Cells.Item(6, 19).Value = 0.007847307243945376
$ws.Cells.Item(6, 20).Value = 0.007847307243945376

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Apoe"
$ws.Cells.Item(7, 3).Value = "Scarb1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 45.524413
$ws.Cells.Item(7, 8).Value = 136.573239
$ws.Cells.Item(7, 9).Value = 0.009939644832300594
$ws.Cells.Item(7, 10).Value = 0.009939644832300592
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.8713403333333334
$ws.Cells.Item(7, 14).Value = 2.614021
$ws.Cells.Item(7, 15).Value = 0.008339794601633706
$ws.Cells.Item(7, 16).Value = 0.008339794601633706
$ws.Cells.Item(7, 17).Value = 39.66725719822434
$ws.Cells.Item(7, 18).Value = 357.005314784019
$ws.Cells.Item(7, 19).Value = 0.00008289459631457686
$ws.Cells.Item(7, 20).Value = 0.00008289459631457684

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Apoe"
$ws.Cells.Item(8, 3).Value = "Scarb1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 45.524413
$ws.Cells.Item(8, 8).Value = 136.573239
$ws.Cells.Item(8, 9).Value = 0.009939644832300594
$ws.Cells.Item(8, 10).Value = 0.009939644832300592
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 17.88507033333333
$ws.Cells.Item(8, 14).Value = 53.65521099999999
$ws.Cells.Item(8, 15).Value = 0.1711820368112258
$ws.Cells.Item(8, 16).Value = 0.1711820368112258
$ws.Cells.Item(8, 17).Value = 814.2073283887142
$ws.Cells.Item(8, 18).Value = 7327.865955498428
$ws.Cells.Item(8, 19).Value = 0.001701488647573391
$ws.Cells.Item(8, 20).Value = 0.001701488647573391

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Apoe"
$ws.Cells.Item(9, 3).Value = "Scarb1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 45.524413
$ws.Cells.Item(9, 8).Value = 136.573239
$ws.Cells.Item(9, 9).Value = 0.009939644832300594
$ws.Cells.Item(9, 10).Value = 0.009939644832300592
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.237038999999999
$ws.Cells.Item(9, 14).Value = 9.711116999999998
$ws.Cells.Item(9, 15).Value = 0.03098242941905719
$ws.Cells.Item(9, 16).Value = 0.03098242941905719
$ws.Cells.Item(9, 17).Value = 147.364300333107
$ws.Cells.Item(9, 18).Value = 1326.278702997963
$ws.Cells.Item(9, 19).Value = 0.0003079543444672497
$ws.Cells.Item(9, 20).Value = 0.0003079543444672496

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Apoe"
$ws.Cells.Item(10, 3).Value = "Scarb1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4438.215250666667
$ws.Cells.Item(10, 8).Value = 13314.645752
$ws.Cells.Item(10, 9).Value = 0.9690247577915309
$ws.Cells.Item(10, 10).Value = 0.9690247577915307
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 82.48638166666666
$ws.Cells.Item(10, 14).Value = 247.459145
$ws.Cells.Item(10, 15).Value = 0.7894957391680832
$ws.Cells.Item(10, 16).Value = 0.7894957391680832
$ws.Cells.Item(10, 17).Value = 366092.3170853113
$ws.Cells.Item(10, 18).Value = 3294830.853767802
$ws.Cells.Item(10, 19).Value = 0.7650409174247975
$ws.Cells.Item(10, 20).Value = 0.7650409174247974

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Apoe"
$ws.Cells.Item(11, 3).Value = "Scarb1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4438.215250666667
$ws.Cells.Item(11, 8).Value = 13314.645752
$ws.Cells.Item(11, 9).Value = 0.9690247577915309
$ws.Cells.Item(11, 10).Value = 0.9690247577915307
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.8713403333333334
$ws.Cells.Item(11, 14).Value = 2.614021
$ws.Cells.Item(11, 15).Value = 0.008339794601633706
$ws.Cells.Item(11, 16).Value = 0.008339794601633706
$ws.Cells.Item(11, 17).Value = 3867.195955920978
$ws.Cells.Item(11, 18).Value = 34804.7636032888
$ws.Cells.Item(11, 19).Value = 0.008081467443879218
$ws.Cells.Item(11, 20).Value = 0.008081467443879217

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Apoe"
$ws.Cells.Item(12, 3).Value = "Scarb1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4438.215250666667
$ws.Cells.Item(12, 8).Value = 13314.645752
$ws.Cells.Item(12, 9).Value = 0.9690247577915309
$ws.Cells.Item(12, 10).Value = 0.9690247577915307
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 17.88507033333333
$ws.Cells.Item(12, 14).Value = 53.65521099999999
$ws.Cells.Item(12, 15).Value = 0.1711820368112258
$ws.Cells.Item(12, 16).Value = 0.1711820368112258
$ws.Cells.Item(12, 17).Value = 79377.79191264596
$ws.Cells.Item(12, 18).Value = 714400.1272138136
$ws.Cells.Item(12, 19).Value = 0.165879631759259
$ws.Cells.Item(12, 20).Value = 0.165879631759259

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Apoe"
$ws.Cells.Item(13, 3).Value = "Scarb1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4438.215250666667
$ws.Cells.Item(13, 8).Value = 13314.645752
$ws.Cells.Item(13, 9).Value = 0.9690247577915309
$ws.Cells.Item(13, 10).Value = 0.9690247577915307
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 3.237038999999999
$ws.Cells.Item(13, 14).Value = 9.711116999999998
$ws.Cells.Item(13, 15).Value = 0.03098242941905719
$ws.Cells.Item(13, 16).Value = 0.03098242941905719
$ws.Cells.Item(13, 17).Value = 14366.67585680277
$ws.Cells.Item(13, 18).Value = 129300.082711225
$ws.Cells.Item(13, 19).Value = 0.03002274116359509
$ws.Cells.Item(13, 20).Value = 0.03002274116359509

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Apoe"
$ws.Cells.Item(14, 3).Value = "Scarb1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 67.02347933333333
$ws.Cells.Item(14, 8).Value = 201.070438
$ws.Cells.Item(14, 9).Value = 0.01463367753909034
$ws.Cells.Item(14, 10).Value = 0.01463367753909034
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 82.48638166666666
$ws.Cells.Item(14, 14).Value = 247.459145
$ws.Cells.Item(14, 15).Value = 0.7894957391680832
$ws.Cells.Item(14, 16).Value = 0.7894957391680832
$ws.Cells.Item(14, 17).Value = 5528.524296917278
$ws.Cells.Item(14, 18).Value = 49756.71867225551
$ws.Cells.Item(14, 19).Value = 0.0115532260654715
$ws.Cells.Item(14, 20).Value = 0.0115532260654715

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Apoe"
$ws.Cells.Item(15, 3).Value = "Scarb1"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 67.02347933333333
$ws.Cells.Item(15, 8).Value = 201.070438
$ws.Cells.Item(15, 9).Value = 0.01463367753909034
$ws.Cells.Item(15, 10).Value = 0.01463367753909034
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.8713403333333334
$ws.Cells.Item(15, 14).Value = 2.614021
$ws.Cells.Item(15, 15).Value = 0.008339794601633706
$ws.Cells.Item(15, 16).Value = 0.008339794601633706
$ws.Cells.Item(15, 17).Value = 58.40026082346644
$ws.Cells.Item(15, 18).Value = 525.602347411198
$ws.Cells.Item(15, 19).Value = 0.000122041864942554
$ws.Cells.Item(15, 20).Value = 0.000122041864942554

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Apoe"
$ws.Cells.Item(16, 3).Value = "Scarb1"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 67.02347933333333
$ws.Cells.Item(16, 8).Value = 201.070438
$ws.Cells.Item(16, 9).Value = 0.01463367753909034
$ws.Cells.Item(16, 10).Value = 0.01463367753909034
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 17.88507033333333
$ws.Cells.Item(16, 14).Value = 53.65521099999999
$ws.Cells.Item(16, 15).Value = 0.1711820368112258
$ws.Cells.Item(16, 16).Value = 0.1711820368112258
$ws.Cells.Item(16, 17).Value = 1198.71964186138
$ws.Cells.Item(16, 18).Value = 10788.47677675242
$ws.Cells.Item(16, 19).Value = 0.002505022727180171
$ws.Cells.Item(16, 20).Value = 0.002505022727180172

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Apoe"
$ws.Cells.Item(17, 3).Value = "Scarb1"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 67.02347933333333
$ws.Cells.Item(17, 8).Value = 201.070438
$ws.Cells.Item(17, 9).Value = 0.01463367753909034
$ws.Cells.Item(17, 10).Value = 0.01463367753909034
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 3.237038999999999
$ws.Cells.Item(17, 14).Value = 9.711116999999998
$ws.Cells.Item(17, 15).Value = 0.03098242941905719
$ws.Cells.Item(17, 16).Value = 0.03098242941905719
$ws.Cells.Item(17, 17).Value = 216.9576165176939
$ws.Cells.Item(17, 18).Value = 1952.618548659246
$ws.Cells.Item(17, 19).Value = 0.0004533868814961089
$ws.Cells.Item(17, 20).Value = 0.0004533868814961089
